$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 124 (which currently holds "view_ingredient"),
# shifting the existing view_* rows down by one.
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new stored procedure entry.
$ws.Cells.Item(124, 1).Value = "stored procedure"
$ws.Cells.Item(124, 2).Value = "new_order_demand_prediction"
$ws.Cells.Item(124, 3).Value = "procedure that get new order demand prediction"

# Match final selection/view state from the authored edit.
$ws.Range("C124").Select()
